$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New category rows to append after the last existing data row (165).
$newRows = @(
    @{ A = "비거주 복지시설 운영업"; B = "Uncategorized" },
    @{ A = "자동차 차체나 트레일러 제조업"; B = "Uncategorized" },
    @{ A = "개인 및 가정용품 수리업"; B = "Uncategorized" }
)

$startRow = 166
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i

    # Clone the formatting of the last existing row (column A uses the bold /
    # bordered / centered style) onto the new row so the new cell picks up the
    # same style index instead of Excel fabricating a brand-new one.
    $ws.Range("A165").Copy($ws.Range("A$r"))

    $ws.Cells.Item($r, 1).Value = $newRows[$i].A
    $ws.Cells.Item($r, 2).Value = $newRows[$i].B
}

Write-Output "Appended $($newRows.Length) rows starting at row $startRow"
